# Update the GLMM individual-effects results sheet:
#  - refresh fixed-effect estimates (rows 2-6)
#  - refresh the random-effects intercept SD row (row 7, column E only)
#  - append the new random-slope SDs and correlation terms (rows 8-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: fixed / cond / (Intercept) ----
$ws.Range("E2").Value = -1.92127545234655
$ws.Range("F2").Value = 0.0397224822223558
$ws.Range("G2").Value = -48.3674570383534
$ws.Range("H2").Value = 0

# ---- Row 3: fixed / cond / habitat_typeExposed/Low SAV ----
$ws.Range("E3").Value = 0.438129514280025
$ws.Range("F3").Value = 0.076409877089608
$ws.Range("G3").Value = 5.73393821542493
$ws.Range("H3").Value = 0.00000000981250882300578

# ---- Row 4: fixed / cond / habitat_typeMod/Dense SAV ----
$ws.Range("E4").Value = 0.0203086113771808
$ws.Range("F4").Value = 0.016508978404945
$ws.Range("G4").Value = 1.23015554803184
$ws.Range("H4").Value = 0.218638862221832

# ---- Row 5: fixed / cond / habitat_typeShallow/Dense SAV ----
$ws.Range("E5").Value = -0.0321112415839058
$ws.Range("F5").Value = 0.065672982351393
$ws.Range("G5").Value = -0.488956652099183
$ws.Range("H5").Value = 0.624872386245141

# ---- Row 6: fixed / cond / habitat_typeShallow/Low SAV ----
$ws.Range("E6").Value = 0.0615415927748369
$ws.Range("F6").Value = 0.0177482297446012
$ws.Range("G6").Value = 3.46747780823365
$ws.Range("H6").Value = 0.000525367119779436

# ---- Row 7: ran_pars / cond / animal_id / sd__(Intercept) ----
$ws.Range("E7").Value = 0.129632913149002

# ---- New rows 8-22: ran_pars / cond / animal_id.1 / sd__ & cor__ terms ----
$ws.Range("A8").Value = "ran_pars"
$ws.Range("B8").Value = "cond"
$ws.Range("C8").Value = "animal_id.1"
$ws.Range("D8").Value = "sd__habitat_typeDeep/Low SAV"
$ws.Range("E8").Value = 0.110026962608003

$ws.Range("A9").Value = "ran_pars"
$ws.Range("B9").Value = "cond"
$ws.Range("C9").Value = "animal_id.1"
$ws.Range("D9").Value = "sd__habitat_typeExposed/Low SAV"
$ws.Range("E9").Value = 0.110026962608003

$ws.Range("A10").Value = "ran_pars"
$ws.Range("B10").Value = "cond"
$ws.Range("C10").Value = "animal_id.1"
$ws.Range("D10").Value = "sd__habitat_typeMod/Dense SAV"
$ws.Range("E10").Value = 0.110026962608003

$ws.Range("A11").Value = "ran_pars"
$ws.Range("B11").Value = "cond"
$ws.Range("C11").Value = "animal_id.1"
$ws.Range("D11").Value = "sd__habitat_typeShallow/Dense SAV"
$ws.Range("E11").Value = 0.110026962608003

$ws.Range("A12").Value = "ran_pars"
$ws.Range("B12").Value = "cond"
$ws.Range("C12").Value = "animal_id.1"
$ws.Range("D12").Value = "sd__habitat_typeShallow/Low SAV"
$ws.Range("E12").Value = 0.110026962608003

$ws.Range("A13").Value = "ran_pars"
$ws.Range("B13").Value = "cond"
$ws.Range("C13").Value = "animal_id.1"
$ws.Range("D13").Value = "cor__habitat_typeDeep/Low SAV.habitat_typeExposed/Low SAV"
$ws.Range("E13").Value = -0.963051236103355

$ws.Range("A14").Value = "ran_pars"
$ws.Range("B14").Value = "cond"
$ws.Range("C14").Value = "animal_id.1"
$ws.Range("D14").Value = "cor__habitat_typeDeep/Low SAV.habitat_typeMod/Dense SAV"
$ws.Range("E14").Value = 0.9274676833602

$ws.Range("A15").Value = "ran_pars"
$ws.Range("B15").Value = "cond"
$ws.Range("C15").Value = "animal_id.1"
$ws.Range("D15").Value = "cor__habitat_typeDeep/Low SAV.habitat_typeShallow/Dense SAV"
$ws.Range("E15").Value = -0.893198898905955

$ws.Range("A16").Value = "ran_pars"
$ws.Range("B16").Value = "cond"
$ws.Range("C16").Value = "animal_id.1"
$ws.Range("D16").Value = "cor__habitat_typeDeep/Low SAV.habitat_typeShallow/Low SAV"
$ws.Range("E16").Value = 0.860196303677535

$ws.Range("A17").Value = "ran_pars"
$ws.Range("B17").Value = "cond"
$ws.Range("C17").Value = "animal_id.1"
$ws.Range("D17").Value = "cor__habitat_typeExposed/Low SAV.habitat_typeMod/Dense SAV"
$ws.Range("E17").Value = -0.963051236103355

$ws.Range("A18").Value = "ran_pars"
$ws.Range("B18").Value = "cond"
$ws.Range("C18").Value = "animal_id.1"
$ws.Range("D18").Value = "cor__habitat_typeExposed/Low SAV.habitat_typeShallow/Dense SAV"
$ws.Range("E18").Value = 0.9274676833602

$ws.Range("A19").Value = "ran_pars"
$ws.Range("B19").Value = "cond"
$ws.Range("C19").Value = "animal_id.1"
$ws.Range("D19").Value = "cor__habitat_typeExposed/Low SAV.habitat_typeShallow/Low SAV"
$ws.Range("E19").Value = -0.893198898905955

$ws.Range("A20").Value = "ran_pars"
$ws.Range("B20").Value = "cond"
$ws.Range("C20").Value = "animal_id.1"
$ws.Range("D20").Value = "cor__habitat_typeMod/Dense SAV.habitat_typeShallow/Dense SAV"
$ws.Range("E20").Value = -0.963051236103355

$ws.Range("A21").Value = "ran_pars"
$ws.Range("B21").Value = "cond"
$ws.Range("C21").Value = "animal_id.1"
$ws.Range("D21").Value = "cor__habitat_typeMod/Dense SAV.habitat_typeShallow/Low SAV"
$ws.Range("E21").Value = 0.9274676833602

$ws.Range("A22").Value = "ran_pars"
$ws.Range("B22").Value = "cond"
$ws.Range("C22").Value = "animal_id.1"
$ws.Range("D22").Value = "cor__habitat_typeShallow/Dense SAV.habitat_typeShallow/Low SAV"
$ws.Range("E22").Value = -0.963051236103355
